$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update report title (A1): report date changed 2025-09-04 -> 2025-09-05 ---
$ws.Range("A1").Value = "萊爾富 工作統計表  篩選月份：202509   (  製表日期:2025-09-05  )"

# --- Row 37: report description / work-content cells now wrap text ---
$ws.Range("P37").WrapText = $true
$ws.Range("AC37").WrapText = $true

# --- Append new row 38 with the new work ticket ---
$ws.Range("A38").Value = 36
$ws.Range("B38").Value = "維修"
$ws.Range("C38").Value = 2025090857
# D38 looks like a number but must stay text - force a text number format
# before assigning it so it is written out as a shared string, not a number.
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "14145114090401"
$ws.Range("E38").Value = "一般件"
$ws.Range("F38").Value = 4145
$ws.Range("G38").Value = "板橋僑興店"
$ws.Range("H38").Value = "新北市板橋區"
$ws.Range("I38").Value = "2025-09-04 17:09:00"
$ws.Range("J38").Value = "星期四"
$ws.Range("K38").Value = "下午"
$ws.Range("L38").Value = "HL23"
$ws.Range("M38").Value = "HL-TM主機"
$ws.Range("N38").Value = 2302
$ws.Range("O38").Value = "客顯示器畫面不正常"
$ws.Range("P38").Value = "門市告知今日工程師到店協助更換TM2主機(TCX800)後客顯就黑畫面，門市表示查看客顯畫面線路未插，將線路插上後客顯畫面也呈現收銀的畫面，客服開啟觸控校正程式後門市告知兩個畫面都跳成客顯畫面，點選螢幕對應也無反應，協助重啟TM仍異常，VNC查看無出現客顯畫面僅有收銀畫面...請台芝到店協助`n※案14145114090301，台芝回覆:09/04 13:35 更換TCx800主機"
$ws.Range("Q38").Value = "THILF04145"
$ws.Range("R38").Value = "新北一"
$ws.Range("S38").Value = "狄澤洋"
$ws.Range("T38").Value = 1
$ws.Range("U38").Value = "已完工"
$ws.Range("V38").Value = "2025-09-04 17:28:52"
$ws.Range("W38").Value = "2025-09-05 11:08:00"
$ws.Range("X38").Value = "2025-09-05 11:34:00"
$ws.Range("Y38").Value = "2025-09-05 21:28:00"
$ws.Range("Z38").Value = 0.4
$ws.Range("AB38").Value = "到場處理"
$ws.Range("AC38").Value = "線路重插後測試正常"
$ws.Range("AK38").Value = "O"

# Match the row-4 formatting pattern (alternating "white" banding row) for row 38
$ws.Range("A4:AK4").Copy()
$ws.Range("A38:AK38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The long note in P38 contains a line break, which triggers automatic
# row-height growth; reset back to the sheet's default (non-custom) height.
$ws.Rows.Item(38).EntireRow.AutoFit()

# --- Update the print area to include the new row ---
$ws.PageSetup.PrintArea = '$A$1:$AK$38'

# --- Restore gridlines / row-col headers (the engine resets these on save) ---
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.DisplayHeadings = $true

# --- Update selection to match the author's final cursor position ---
$ws.Range("AC35").Select()
